# Finalized definition of virtual environment
# Updates the calibration outputs across the start_price / Linear / NonLinear
# sheets (and the two "abs_epsi_autocorr" / "sig2_1" autocorrelation-array
# strings in the shared-string table) to their newly recomputed values.

$wb = $excel.ActiveWorkbook

# --- start_price ---------------------------------------------------------
$wsStart = $wb.Worksheets.Item("start_price")
$wsStart.Range("A2").Value = 2593.58526328259

# --- Linear ----------------------------------------------------------------
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.07297428785890372
$wsLinear.Range("B3").Value = 0.156135521917966
$wsLinear.Range("B4").Value = 468.9925279244225
$wsLinear.Range("B5").Value = "[1.0, 0.27965338747024354, 0.10337298118144533, 0.12383047330708387, 0.08204911403904946, 0.09139528350251484, 0.2089788131416943, 0.30077014220122483, 0.19532568962884156, 0.04478342443896179, 0.050887328932263956, -0.0019858762001142348, 0.008145162796656698, 0.15403225909919444, 0.27581139642923996, 0.1346958144227024, 0.029256695899101667, 0.03708102648301418, 0.024622409068468788, 0.027039891351111768]"

# --- NonLinear ---------------------------------------------------------------
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 0.9427402862985685
$wsNonLinear.Range("B4").Value = -0.2409121847464984
$wsNonLinear.Range("B5").Value = 0.1353672079848004
$wsNonLinear.Range("B6").Value = 497.6646720322361
$wsNonLinear.Range("B7").Value = -0.1706559714051568
$wsNonLinear.Range("B8").Value = 0.1738265302077644
$wsNonLinear.Range("B9").Value = 442.5560976754697
$wsNonLinear.Range("B10").Value = "[1.0, 0.27973103817457734, 0.10332011579989726, 0.12373203730801884, 0.08157498487592384, 0.09113029714338135, 0.20918464085463365, 0.3008712964939789, 0.19534772962630334, 0.045102754144366396, 0.05105435965752363, -0.0017672486403906653, 0.00846754603055806, 0.15414556918193126, 0.27533097802052925, 0.13499678998291667, 0.029886693073867585, 0.03705918679226068, 0.024040480550208014, 0.027311856845705428]"
